$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 531, shifting existing rows 531..619 down to 532..620.
$ws.Rows(531).Insert()

# Populate the newly inserted row 531 with a fresh weekly price entry
# (same categorical/meta columns as the row that used to sit at 531,
# but new date + price figures).
$ws.Range("A531").Value = 3
$ws.Range("B531").Value = "Femacal de La Calera"
$ws.Range("C531").Value = "Coquimbo"
$ws.Range("D531").Value = "2023-09-11"
$ws.Range("E531").Value = 5
$ws.Range("F531").Value = 100112009
$ws.Range("G531").Value = "Acelga"
$ws.Range("H531").Value = "Sin especificar"
$ws.Range("I531").Value = "Primera"
$ws.Range("J531").Value = 230
$ws.Range("K531").Value = 4000
$ws.Range("L531").Value = 4500
$ws.Range("M531").Value = 4239
$ws.Range("N531").Value = "$/docena de atados (6 kilos)"
$ws.Range("O531").Value = "Provincia de Quillota"
$ws.Range("P531").Value = 706
$ws.Range("Q531").Value = 6
$ws.Range("R531").Value = "Hortaliza"
